# Update "想去人数" (column F) figures on both the "展览" (Exhibition) sheet
# and the "全部类型" (All Types) sheet to reflect newly scraped attendance counts.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 53
$wsExhibition.Range("F3").Value = 779
$wsExhibition.Range("F4").Value = 38
$wsExhibition.Range("F6").Value = 63
$wsExhibition.Range("F7").Value = 270
$wsExhibition.Range("F8").Value = 3858
$wsExhibition.Range("F9").Value = 86
$wsExhibition.Range("F10").Value = 4547
$wsExhibition.Range("F11").Value = 494
$wsExhibition.Range("F12").Value = 1146
$wsExhibition.Range("F13").Value = 69

# --- Sheet "全部类型" (All Types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 53
$wsAll.Range("F3").Value = 779
$wsAll.Range("F4").Value = 38
$wsAll.Range("F6").Value = 63
$wsAll.Range("F8").Value = 270
$wsAll.Range("F9").Value = 3858
$wsAll.Range("F10").Value = 86
$wsAll.Range("F11").Value = 4547
$wsAll.Range("F12").Value = 494
$wsAll.Range("F13").Value = 1146
$wsAll.Range("F14").Value = 69

$wb.Save()
